$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.131.44"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").Value = "2.666.33"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.41"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.18"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").Value = "2.667.25"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.80"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "3.164.85"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").Value = "67.235.23"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "2.670.29"
$ws.Range("E18").Value = "  -3.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.64"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "361.64"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.51"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.80"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("E24").Value = "  -4.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.47"
$ws.Range("E25").Value = "  -4.90%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "550.66"
$ws.Range("E31").Value = "  -5.33%  "
$ws.Range("E32").Value = "  -3.85%  "
$ws.Range("E33").Value = "  -3.31%  "
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  -4.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.54"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "155.69"
$ws.Range("E39").Value = "  -3.89%  "
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.83"
$ws.Range("E41").Value = "  -3.33%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.27"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.90"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("E44").Value = "  -4.83%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.26"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").Value = "0.0₆0298"
$ws.Range("E47").Value = "  -6.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.585"
$ws.Range("E48").Value = "  -3.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.22"
$ws.Range("E49").Value = "  -4.07%  "
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("E51").Value = "  -2.93%  "
